# Applies the cryptos.xlsx price/volume refresh described by the commit
# "Updated cryptos list on Sat Mar 30 05:29:26 UTC 2024 with GitHub Actions".
#
# All target cells hold plain text (prices/percentages formatted as text,
# e.g. "1.00", "70.005.13", "  -0.55%  "). A naive `.Value = "..."` assignment
# gets silently re-interpreted as a number by Excel (losing exact formatting,
# e.g. "607.20" -> 607.2000000000001), so instead we set a self-quoting
# formula (="text") and immediately flatten it to a static value via
# Copy + PasteSpecial(xlPasteValues), which preserves the text exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-ExactText($addr, $val) {
    $r = $ws.Range($addr)
    $r.Formula = '="' + $val.Replace('"', '""') + '"'
    $r.Copy()
    $r.PasteSpecial(-4163)  # xlPasteValues
}

Set-ExactText 'D2' '69.940.55'
Set-ExactText 'E2' '  -0.42%  '
Set-ExactText 'D3' '3.516.41'
Set-ExactText 'E3' '  -1.26%  '
Set-ExactText 'D4' '0.999'
Set-ExactText 'E4' '  -0.13%  '
Set-ExactText 'D5' '607.20'
Set-ExactText 'E5' '  -0.10%  '
Set-ExactText 'D6' '198.65'
Set-ExactText 'E6' '  +6.47%  '
Set-ExactText 'D7' '0.627'
Set-ExactText 'E7' '  +1.14%  '
Set-ExactText 'E8' '  -0.09%  '
Set-ExactText 'D9' '0.214'
Set-ExactText 'E9' '  -0.50%  '
Set-ExactText 'D10' '0.659'
Set-ExactText 'E10' '  +2.01%  '
Set-ExactText 'D11' '54.26'
Set-ExactText 'E11' '  +0.58%  '
Set-ExactText 'D12' '0.0000306'
Set-ExactText 'E12' '  -0.90%  '
Set-ExactText 'D13' '9.65'
Set-ExactText 'E13' '  +2.08%  '
Set-ExactText 'D14' '4.074.25'
Set-ExactText 'E14' '  -1.34%  '
Set-ExactText 'D15' '598.59'
Set-ExactText 'E15' '  +4.79%  '
Set-ExactText 'D16' '70.091.99'
Set-ExactText 'E16' '  -0.34%  '
Set-ExactText 'D17' '19.05'
Set-ExactText 'E17' '  +0.53%  '
Set-ExactText 'D18' '12.69'
Set-ExactText 'E18' '  +0.13%  '
Set-ExactText 'D19' '3.513.61'
Set-ExactText 'E19' '  -1.96%  '
Set-ExactText 'E20' '  -0.19%  '
Set-ExactText 'D21' '1.00'
Set-ExactText 'E21' '  +0.91%  '
Set-ExactText 'D22' '17.78'
Set-ExactText 'E22' '  +1.98%  '
Set-ExactText 'D23' '103.91'
Set-ExactText 'E23' '  +10.83%  '
Set-ExactText 'E24' '  -1.90%  '
Set-ExactText 'D25' '5.07'
Set-ExactText 'E25' '  +3.73%  '
Set-ExactText 'D26' '3.12'
Set-ExactText 'E26' '  +5.96%  '
Set-ExactText 'D27' '11.03'
Set-ExactText 'E27' '  +1.00%  '
Set-ExactText 'D28' '9.83'
Set-ExactText 'E28' '  +4.99%  '
Set-ExactText 'D29' '33.91'
Set-ExactText 'E29' '  +4.93%  '
Set-ExactText 'D30' '4.52'
Set-ExactText 'E30' '  +21.37%  '
Set-ExactText 'D31' '7.23'
Set-ExactText 'E31' '  +2.54%  '
Set-ExactText 'D32' '12.77'
Set-ExactText 'E32' '  +4.56%  '
Set-ExactText 'E33' '  +1.49%  '
Set-ExactText 'D34' '63.92'
Set-ExactText 'E34' '  -0.28%  '
Set-ExactText 'D35' '3.722.62'
Set-ExactText 'E35' '  +2.35%  '
Set-ExactText 'D36' '522.13'
Set-ExactText 'E36' '  -0.51%  '
Set-ExactText 'E37' '  -0.03%  '
Set-ExactText 'D38' '0.0₃0801'
Set-ExactText 'E38' '  +2.37%  '
Set-ExactText 'D39' '3.03'
Set-ExactText 'E39' '  -4.89%  '
Set-ExactText 'B40' 'TheGraph'
Set-ExactText 'C40' 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-ExactText 'D40' '0.393'
Set-ExactText 'E40' '  -2.62%  '
Set-ExactText 'B41' 'InjectiveProtocol'
Set-ExactText 'C41' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-ExactText 'D41' '37.01'
Set-ExactText 'E41' '  -1.21%  '
Set-ExactText 'E42' '  +1.59%  '
Set-ExactText 'E43' '  -0.46%  '
Set-ExactText 'D44' '0.0464'
Set-ExactText 'E44' '  +1.69%  '
Set-ExactText 'D45' '2.88'
Set-ExactText 'E45' '  -2.43%  '
Set-ExactText 'E46' '  +0.99%  '
Set-ExactText 'E47' '  -6.08%  '
Set-ExactText 'E48' '  -4.13%  '
Set-ExactText 'E49' '  +0.18%  '
Set-ExactText 'D50' '132.54'
Set-ExactText 'E50' '  -2.03%  '
Set-ExactText 'B51' 'FLOKI'
Set-ExactText 'C51' 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
Set-ExactText 'D51' '0.000240'
Set-ExactText 'E51' '  -2.19%  '
